# Utilities.xlsx - "clarify status report for unit tests"
#
# The unit-test row on Sheet1 (A3:E3) compares an "Expected" qlVersion()
# string (B3) against the "Actual" qlVersion() returned by the QuantLib
# add-in (C3, via the _xll.qlVersion() UDF cached in E3), and reports
# PASS/FAIL in D3.
#
# This edit:
#   1. Updates the "Expected" version in B3 to "1.8.1" (previously it was
#      pinned to the same string as the actual result, "1.8.2", so the
#      test always passed trivially).
#   2. Clarifies the status formula in D3 so a broken/erroring lookup is
#      reported as "ERROR" instead of being silently folded into "FAIL".
#
# Note: E3's own formula (=_xll.qlVersion()) is left exactly as-is - it
# calls straight into the add-in and is not something this script needs
# to (or safely can) re-derive.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. "Expected" version bumped down to 1.8.1.
$ws.Range("B3").Value = "1.8.1"

# 2. Status formula: distinguish a hard ERROR (either side errored out)
#    from a plain FAIL (values just don't match).
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
